$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 92.13
$ws.Range("F2").Value = 4.69
$ws.Range("K2").Value = 62.2
$ws.Range("N2").Value = 85.96878041621773

$ws.Range("K3").Value = 55.8
$ws.Range("N3").Value = 85.96878041621773
